$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Addr, $Val)
    $cell = $Sheet.Range($Addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = "'" + $Val
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '26.165.93'
Set-TextValue $ws 'E2' '  -0.07%  '
Set-TextValue $ws 'D3' '1.655.82'
Set-TextValue $ws 'E3' '  -0.20%  '
Set-TextValue $ws 'E4' '  -0.26%  '
Set-TextValue $ws 'D5' '218.65'
Set-TextValue $ws 'E5' '  -0.38%  '
Set-TextValue $ws 'D6' '0.5236'
Set-TextValue $ws 'E6' '  +0.08%  '
Set-TextValue $ws 'D7' '1.004'
Set-TextValue $ws 'E7' '  -0.25%  '
Set-TextValue $ws 'D8' '0.2664'
Set-TextValue $ws 'E8' '  +1.31%  '
Set-TextValue $ws 'D9' '0.06345'
Set-TextValue $ws 'E9' '  +0.75%  '
Set-TextValue $ws 'D10' '20.54'
Set-TextValue $ws 'E10' '  -0.38%  '
Set-TextValue $ws 'D11' '0.07685'
Set-TextValue $ws 'E11' '  -1.73%  '
Set-TextValue $ws 'E12' '  +2.93%  '
Set-TextValue $ws 'D13' '1.704.58'
Set-TextValue $ws 'E13' '  +2.87%  '
Set-TextValue $ws 'D14' '1.885.03'
Set-TextValue $ws 'E14' '  -0.10%  '
Set-TextValue $ws 'D15' '0.5611'
Set-TextValue $ws 'E15' '  +1.14%  '
Set-TextValue $ws 'D16' '0.0₅8183'
Set-TextValue $ws 'E16' '  +2.08%  '
Set-TextValue $ws 'E17' '  +0.60%  '
Set-TextValue $ws 'D18' '26.156.65'
Set-TextValue $ws 'E19' '  -0.23%  '
Set-TextValue $ws 'D20' '4.655'
Set-TextValue $ws 'E20' '  +0.36%  '
Set-TextValue $ws 'E21' '  +3.64%  '
Set-TextValue $ws 'D22' '192.59'
Set-TextValue $ws 'E22' '  -1.78%  '
Set-TextValue $ws 'D23' '5.958'
Set-TextValue $ws 'E23' '  +0.07%  '
Set-TextValue $ws 'E24' '  -0.30%  '
Set-TextValue $ws 'D25' '145.44'
Set-TextValue $ws 'E25' '  -0.60%  '
Set-TextValue $ws 'D26' '0.1194'
Set-TextValue $ws 'E26' '  -0.56%  '
Set-TextValue $ws 'D27' '7.267'
Set-TextValue $ws 'E27' '  +1.73%  '
Set-TextValue $ws 'D28' '15.96'
Set-TextValue $ws 'E28' '  -0.37%  '
Set-TextValue $ws 'D29' '1.519'
Set-TextValue $ws 'E29' '  +1.87%  '
Set-TextValue $ws 'D30' '0.05479'
Set-TextValue $ws 'E30' '  -4.58%  '
Set-TextValue $ws 'E31' '  -0.33%  '
Set-TextValue $ws 'E32' '  -0.54%  '
Set-TextValue $ws 'D33' '3.368'
Set-TextValue $ws 'E33' '  +0.01%  '
Set-TextValue $ws 'D34' '1.563'
Set-TextValue $ws 'E34' '  -1.38%  '
Set-TextValue $ws 'D35' '0.9500'
Set-TextValue $ws 'E35' '  -0.58%  '
Set-TextValue $ws 'D36' '2.779'
Set-TextValue $ws 'E37' '  -0.75%  '
Set-TextValue $ws 'D38' '0.5694'
Set-TextValue $ws 'E38' '  -0.42%  '
Set-TextValue $ws 'D39' '0.01591'
Set-TextValue $ws 'E39' '  -0.38%  '
Set-TextValue $ws 'D40' '5.871'
Set-TextValue $ws 'E40' '  -1.54%  '
Set-TextValue $ws 'D41' '1.003'
Set-TextValue $ws 'E41' '  -0.21%  '
Set-TextValue $ws 'B42' 'Maker'
Set-TextValue $ws 'C42' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws 'D42' '1.031.05'
Set-TextValue $ws 'E42' '  -3.12%  '
Set-TextValue $ws 'B43' 'TrustWalletToken'
Set-TextValue $ws 'C43' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws 'D43' '0.8334'
Set-TextValue $ws 'E43' '  -1.70%  '
Set-TextValue $ws 'D44' '100.97'
Set-TextValue $ws 'E44' '  -2.91%  '
Set-TextValue $ws 'D45' '1.794.82'
Set-TextValue $ws 'E45' '  -0.17%  '
Set-TextValue $ws 'D46' '58.00'
Set-TextValue $ws 'E46' '  -0.27%  '
Set-TextValue $ws 'B47' 'Frax'
Set-TextValue $ws 'C47' 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws 'D47' '0.9978'
Set-TextValue $ws 'E47' '  -0.94%  '
Set-TextValue $ws 'B48' 'EnergySwap'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D48' '8.048'
Set-TextValue $ws 'E48' '  +0.37%  '
Set-TextValue $ws 'B49' 'Mantle'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws 'D49' '0.4346'
Set-TextValue $ws 'E49' '  -1.35%  '
Set-TextValue $ws 'B50' 'BabyDogeCoin'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws 'D50' '0.0₈103'
Set-TextValue $ws 'E50' '  +0.30%  '
Set-TextValue $ws 'D51' '0.05214'
Set-TextValue $ws 'E51' '  +0.23%  '
